# Insert two new rows above row 61 (before "econ_start_time") on the
# "constants" sheet to add the new plot_start_time / plot_end_time
# parameters, as part of tidying the main outputs plotting function.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert two rows at row 61; shifts old row 61 (econ_start_time) and
# everything below it down to row 63 onward.
$ws.Rows.Item(61).Resize(2).Insert()

# Only columns A:E are used on this sheet - copy formatting from the
# range directly below (old row 61, now row 63, which shares the exact
# same formatting as row 60 above the insertion point) onto just those
# five-wide ranges instead of doing a whole-row copy/paste (which would
# stamp formatting all the way out to column XFD).
$ws.Range("A63:E63").Copy()
$ws.Range("A61:E61").PasteSpecial(-4122) | Out-Null
$ws.Range("A63:E63").Copy()
$ws.Range("A62:E62").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new shared strings in the same order they appear in the
# saved workbook's sharedStrings table (plot_end_time/description pair
# first, then plot_start_time/description pair), then write the values
# into their actual cells (row 61 = plot_start_time, row 62 = plot_end_time).
$ws.Cells.Item(1, 10).Value = "plot_end_time"
$ws.Cells.Item(1, 10).Value = ""
$ws.Cells.Item(2, 10).Value = "Time that plots finish at"
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(3, 10).Value = "plot_start_time"
$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(4, 10).Value = "Time that plots start from"
$ws.Cells.Item(4, 10).Value = ""

# New row 61: plot_start_time
$ws.Cells.Item(61, 1).Value = "plot_start_time"
$ws.Cells.Item(61, 2).Value = 1990
$ws.Cells.Item(61, 5).Value = "Time that plots start from"

# New row 62: plot_end_time
$ws.Cells.Item(62, 1).Value = "plot_end_time"
$ws.Cells.Item(62, 2).Value = 2035
$ws.Cells.Item(62, 5).Value = "Time that plots finish at"

# Match the scrolled viewport / selection recorded in the saved workbook.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("A61").Select()
